$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Story "I wan to send bills to other housemates" (row 9) progressed to IN PROGRESS,
# and "I want to view monthly bill" (row 8) is now DONE, reflecting that the
# "send bill" button/feature has been wired up.

# Row 9: NOT STARTED -> IN PROGRESS. Copy the current (pre-edit) formatting of C8,
# which already carries the "IN PROGRESS" look (fill/font/border), onto C9.
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C9").Value = "IN PROGRESS"

# Row 8: IN PROGRESS -> DONE. Copy the "DONE" formatting from one of the
# already-DONE status cells (C2) onto C8.
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C8").Value = "DONE"

$excel.CutCopyMode = $false

# Reflect the updated active selection seen in the saved workbook.
$ws.Range("F12").Select()
